$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(6, 9).Value = 'aa'
$ws.Cells.Item(6, 10).Value = 'Agree/Accept'
$ws.Cells.Item(34, 9).Value = 'sd'
$ws.Cells.Item(34, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(36, 9).Value = 'sd'
$ws.Cells.Item(36, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(38, 9).Value = 'sv'
$ws.Cells.Item(38, 10).Value = 'Statement-opinion'
$ws.Cells.Item(42, 9).Value = 'sv'
$ws.Cells.Item(42, 10).Value = 'Statement-opinion'
$ws.Cells.Item(44, 9).Value = 'sv'
$ws.Cells.Item(44, 10).Value = 'Statement-opinion'
$ws.Cells.Item(46, 9).Value = 'sd'
$ws.Cells.Item(46, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(51, 9).Value = 'aa'
$ws.Cells.Item(51, 10).Value = 'Agree/Accept'
$ws.Cells.Item(53, 9).Value = 'sd'
$ws.Cells.Item(53, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(57, 9).Value = 'aa'
$ws.Cells.Item(57, 10).Value = 'Agree/Accept'
$ws.Cells.Item(58, 9).Value = 'sv'
$ws.Cells.Item(58, 10).Value = 'Statement-opinion'
$ws.Cells.Item(62, 9).Value = 'sd'
$ws.Cells.Item(62, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(63, 9).Value = 'ba'
$ws.Cells.Item(63, 10).Value = 'Appreciation'
$ws.Cells.Item(67, 9).Value = 'sv'
$ws.Cells.Item(67, 10).Value = 'Statement-opinion'
$ws.Cells.Item(69, 9).Value = 'aa'
$ws.Cells.Item(69, 10).Value = 'Agree/Accept'
$ws.Cells.Item(73, 9).Value = 'sv'
$ws.Cells.Item(73, 10).Value = 'Statement-opinion'
$ws.Cells.Item(74, 9).Value = 'sv'
$ws.Cells.Item(74, 10).Value = 'Statement-opinion'
$ws.Cells.Item(75, 9).Value = 'sv'
$ws.Cells.Item(75, 10).Value = 'Statement-opinion'
$ws.Cells.Item(78, 9).Value = 'sv'
$ws.Cells.Item(78, 10).Value = 'Statement-opinion'
$ws.Cells.Item(79, 9).Value = 'aa'
$ws.Cells.Item(79, 10).Value = 'Agree/Accept'
$ws.Cells.Item(81, 9).Value = 'ba'
$ws.Cells.Item(81, 10).Value = 'Appreciation'
$ws.Cells.Item(90, 9).Value = 'b'
$ws.Cells.Item(90, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(98, 9).Value = 'sd'
$ws.Cells.Item(98, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(106, 9).Value = 'sd'
$ws.Cells.Item(106, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(108, 9).Value = 'sd'
$ws.Cells.Item(108, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(111, 9).Value = 'sd'
$ws.Cells.Item(111, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(120, 9).Value = 'sd'
$ws.Cells.Item(120, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(128, 9).Value = 'b'
$ws.Cells.Item(128, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(130, 9).Value = 'b'
$ws.Cells.Item(130, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(138, 9).Value = 'sv'
$ws.Cells.Item(138, 10).Value = 'Statement-opinion'
$ws.Cells.Item(147, 9).Value = 'aa'
$ws.Cells.Item(147, 10).Value = 'Agree/Accept'
$ws.Cells.Item(154, 9).Value = 'aa'
$ws.Cells.Item(154, 10).Value = 'Agree/Accept'
